# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 15; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 25; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 44; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 46; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 52; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 58; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 60; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 61; I = 'ba'; J = 'Appreciation' }
    @{ Row = 73; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 78; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 88; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 92; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 98; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 100; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 114; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 116; I = 'qy'; J = 'Yes-No-Question' }
    @{ Row = 120; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 121; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 124; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 127; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 130; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 131; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 134; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 136; I = '%'; J = 'Uninterpretable' }
    @{ Row = 146; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 150; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 152; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 158; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 171; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 174; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 189; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 209; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 235; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 242; I = 'ba'; J = 'Appreciation' }
    @{ Row = 243; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 275; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 278; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 281; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 285; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 301; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 302; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 307; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 315; I = '%'; J = 'Uninterpretable' }
    @{ Row = 326; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 327; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 334; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 340; I = '%'; J = 'Uninterpretable' }
    @{ Row = 349; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 352; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 367; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 373; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 374; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 378; I = 'aa'; J = 'Agree/Accept' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
